$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 42 data rows to 45. Insert a row at the top for the new
# "Date and Time" entry, one row before "Idling time percentage" for the new
# "Cycle Count of battery" entry, and two rows at the bottom for the two newly
# appended speed-bucket entries.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(44).Insert()

# Write the final label/value for every row so the content matches exactly.
$ws.Cells.Item(1, 1).Value = "Date and Time"
$ws.Cells.Item(1, 2).Value = "2024-03-12 10:41:58.064000 to 2024-03-12 11:48:27.259000"

$ws.Cells.Item(2, 1).Value = "Total time taken for the ride"
$ws.Cells.Item(2, 2).Value = 0.04619045138888889

$ws.Cells.Item(3, 1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(3, 2).Value = 32.33816583333333

$ws.Cells.Item(4, 1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(4, 2).Value = 1642.223101016667

$ws.Cells.Item(5, 1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(5, 2).Value = 7.43

$ws.Cells.Item(6, 1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(6, 2).Value = 7.051

$ws.Cells.Item(7, 1).Value = "Starting SoC (%)"
$ws.Cells.Item(7, 2).Value = 17

$ws.Cells.Item(8, 1).Value = "Ending SoC (%)"
$ws.Cells.Item(8, 2).Value = 99

$ws.Cells.Item(9, 1).Value = "Total distance covered (km)"
$ws.Cells.Item(9, 2).Value = 36.21829742812135

$ws.Cells.Item(10, 1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10, 2).Value = 45.34236056445818

$ws.Cells.Item(11, 1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(11, 2).Value = 82

$ws.Cells.Item(12, 1).Value = "Mode"
$ws.Cells.Item(12, 2).Value = "Custom mode`n79.89%`nEco mode`n18.74%`nSports mode`n0.03%"

$ws.Cells.Item(13, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(13, 2).Value = 5233.3208

$ws.Cells.Item(14, 1).Value = "Average Power(kW)"
$ws.Cells.Item(14, 2).Value = -1489.920152132056

$ws.Cells.Item(15, 1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(15, 2).Value = 3.365682143888889

$ws.Cells.Item(16, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(16, 2).Value = 0.2045275331437714

$ws.Cells.Item(17, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 3.491

$ws.Cells.Item(18, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(18, 2).Value = 3.032

$ws.Cells.Item(19, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19, 2).Value = 0.4590000000000001

$ws.Cells.Item(20, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20, 2).Value = 31

$ws.Cells.Item(21, 1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(21, 2).Value = 47

$ws.Cells.Item(22, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(22, 2).Value = 16

$ws.Cells.Item(23, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23, 2).Value = 60

$ws.Cells.Item(24, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24, 2).Value = 63

$ws.Cells.Item(25, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25, 2).Value = 60

$ws.Cells.Item(26, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26, 2).Value = 59

$ws.Cells.Item(27, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27, 2).Value = 98

$ws.Cells.Item(28, 1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(28, 2).Value = 0

$ws.Cells.Item(29, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(29, 2).Value = 47

$ws.Cells.Item(30, 1).Value = "lowest cell temp(C)"
$ws.Cells.Item(30, 2).Value = 31

$ws.Cells.Item(31, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(31, 2).Value = 16

$ws.Cells.Item(32, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(32, 2).Value = 54

$ws.Cells.Item(33, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(33, 2).Value = 1.746260955

$ws.Cells.Item(34, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(34, 2).Value = [double]"1.215720520050125e-07"

$ws.Cells.Item(35, 1).Value = "Cycle Count of battery"
$ws.Cells.Item(35, 2).Value = 53

$ws.Cells.Item(36, 1).Value = "Idling time percentage"
$ws.Cells.Item(36, 2).Value = 24.73596594165325

$ws.Cells.Item(37, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(37, 2).Value = 5.108752012662719

$ws.Cells.Item(38, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(38, 2).Value = 6.153972109270529

$ws.Cells.Item(39, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(39, 2).Value = 9.349671151379527

$ws.Cells.Item(40, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(40, 2).Value = 14.190977812952

$ws.Cells.Item(41, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(41, 2).Value = 8.708348115601888

$ws.Cells.Item(42, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(42, 2).Value = 6.536036896542313

$ws.Cells.Item(43, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(43, 2).Value = 14.48571350599023

$ws.Cells.Item(44, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(44, 2).Value = 10.61867205196081

$ws.Cells.Item(45, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(45, 2).Value = 0
